# Weekly refresh of the Alcachofa (Hortaliza) price records for
# "Terminal La Palmera de La Serena". The underlying data rows are the
# same set of weekly observations, just re-dated/re-ordered as a new
# week's record set rotates into the table (row 5 is unaffected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44484
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9500
$ws.Range("O2").Value = 'Provincia del Elquí'
$ws.Range("P2").Value = 317

# Row 3
$ws.Range("D3").Value = 44420
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("P3").Value = 362

# Row 4
$ws.Range("D4").Value = 44420
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("P4").Value = 338
$ws.Range("Q4").Value = 40

# Row 6
$ws.Range("D6").Value = 44427
$ws.Range("H6").Value = 'Madrigal'
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("N6").Value = '$/caja 40 unidades'
$ws.Range("O6").Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 312
$ws.Range("Q6").Value = 40

# Row 7
$ws.Range("D7").Value = 44426
$ws.Range("H7").Value = 'Española'
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 11500
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11750
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("P7").Value = 392
$ws.Range("Q7").Value = 30

# Row 8
$ws.Range("D8").Value = 44426
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 12500
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12750
$ws.Range("O8").Value = 'Provincia de Limarí'
$ws.Range("P8").Value = 319

# Row 9
$ws.Range("D9").Value = 44438
$ws.Range("H9").Value = 'Española'
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11500
$ws.Range("N9").Value = '$/caja 30 unidades'
$ws.Range("O9").Value = 'Provincia del Elquí'
$ws.Range("P9").Value = 383
$ws.Range("Q9").Value = 30
